$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "277.89"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.76%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.22"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.81%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.874"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.53%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06363"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.61%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.970"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.56%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.251"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-7.24%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8816"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.25%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1523"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.38%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05104"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.91%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07546"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.94%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02966"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.64%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09012"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.33%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001577"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.33%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006429"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.71%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005895"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.81%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.461"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.41%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.315"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.17%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.50%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.78%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.28%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.903"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.72%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04419"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.89%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001172"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.42%"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.16%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.18%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "14.12%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04147"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.66%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006863"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.78%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1179"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.21%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002021"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-10.00%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01120"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-11.20%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005182"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.80%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.487"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.02025"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-4.44%"
